# Scheduled-runner refresh of market-price-derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-class
# leve-profit sheets. Values only; no formulas/formatting involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 219.58333
$ws.Range("I96").Value = 189.55556
$ws.Range("K96").Value = 568.66668
$ws.Range("M96").Value = 804.33332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H68").Value = 37444
$ws.Range("J68").Value = 37444
$ws.Range("L68").Value = 37444
$ws.Range("N68").Value = -39066

$ws.Range("H71").Value = 37444
$ws.Range("J71").Value = 37444
$ws.Range("L71").Value = 112332
$ws.Range("N71").Value = -120444

$ws.Range("H80").Value = 25263.309
$ws.Range("J80").Value = 25263.309
$ws.Range("L80").Value = 25263.309
$ws.Range("N80").Value = -27259.309

$ws.Range("H83").Value = 25263.309
$ws.Range("J83").Value = 25263.309
$ws.Range("L83").Value = 75789.927
$ws.Range("N83").Value = -85773.927

$ws.Range("H97").Value = 2621.8845
$ws.Range("I97").Value = 2429.5
$ws.Range("J97").Value = 3263.1667
$ws.Range("K97").Value = 2429.5
$ws.Range("L97").Value = 3263.1667
$ws.Range("M97").Value = -1933.5
$ws.Range("N97").Value = -4255.1667

$ws.Range("H102").Value = 58825988
$ws.Range("I102").Value = 2573
$ws.Range("J102").Value = 142859440
$ws.Range("K102").Value = 2573
$ws.Range("L102").Value = 142859440
$ws.Range("M102").Value = -951
$ws.Range("N102").Value = -142862684

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 45000
$ws.Range("J62").Value = 45000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46372

$ws.Range("H65").Value = 45000
$ws.Range("J65").Value = 45000
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141864

$ws.Range("H86").Value = 1585.5682
$ws.Range("I86").Value = 1621.0555
$ws.Range("J86").Value = 1425.875
$ws.Range("K86").Value = 1621.0555
$ws.Range("L86").Value = 1425.875
$ws.Range("M86").Value = -498.0554999999999
$ws.Range("N86").Value = -3671.875

$ws.Range("H89").Value = 1585.5682
$ws.Range("I89").Value = 1621.0555
$ws.Range("J89").Value = 1425.875
$ws.Range("K89").Value = 8105.2775
$ws.Range("L89").Value = 7129.375
$ws.Range("M89").Value = -2489.2775
$ws.Range("N89").Value = -18361.375

$ws.Range("H94").Value = 925.5
$ws.Range("I94").Value = 854.2353000000001
$ws.Range("J94").Value = 1098.5714
$ws.Range("K94").Value = 854.2353000000001
$ws.Range("L94").Value = 1098.5714
$ws.Range("M94").Value = -403.2353000000001
$ws.Range("N94").Value = -2000.5714

$ws.Range("H105").Value = 3452.8823
$ws.Range("I105").Value = 1958.8
$ws.Range("J105").Value = 5587.2856
$ws.Range("K105").Value = 1958.8
$ws.Range("L105").Value = 5587.2856
$ws.Range("M105").Value = -211.8
$ws.Range("N105").Value = -9081.285599999999

$ws.Range("H107").Value = 4524.2856
$ws.Range("I107").Value = 4702.794
$ws.Range("K107").Value = 4702.794
$ws.Range("M107").Value = -2782.794

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 212.68182
$ws.Range("I22").Value = 210.15384
$ws.Range("J22").Value = 216.33333
$ws.Range("K22").Value = 210.15384
$ws.Range("L22").Value = 216.33333
$ws.Range("M22").Value = 139.84616
$ws.Range("N22").Value = -916.3333299999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 12694.538
$ws.Range("I9").Value = 201
$ws.Range("J9").Value = 13735.667
$ws.Range("K9").Value = 603
$ws.Range("L9").Value = 41207.001
$ws.Range("M9").Value = -379
$ws.Range("N9").Value = -41655.001

$ws.Range("H12").Value = 671.2174
$ws.Range("I12").Value = 195
$ws.Range("J12").Value = 879.5625
$ws.Range("K12").Value = 585
$ws.Range("L12").Value = 2638.6875
$ws.Range("M12").Value = -412
$ws.Range("N12").Value = -2984.6875

$ws.Range("H98").Value = 451.75
$ws.Range("I98").Value = 203
$ws.Range("J98").Value = 534.6667
$ws.Range("K98").Value = 609
$ws.Range("L98").Value = 1604.0001
$ws.Range("M98").Value = 889
$ws.Range("N98").Value = -4600.0001

$ws.Range("H129").Value = 54024.74
$ws.Range("I129").Value = 84074.164
$ws.Range("J129").Value = 2511.4285
$ws.Range("K129").Value = 252222.492
$ws.Range("L129").Value = 7534.2855
$ws.Range("M129").Value = -247222.492
$ws.Range("N129").Value = -17534.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3221.25
$ws.Range("J80").Value = 3399.4443
$ws.Range("L80").Value = 3399.4443
$ws.Range("N80").Value = -5395.4443

$ws.Range("H83").Value = 3221.25
$ws.Range("J83").Value = 3399.4443
$ws.Range("L83").Value = 16997.2215
$ws.Range("N83").Value = -26981.2215

$ws.Range("H97").Value = 3227.5
$ws.Range("I97").Value = 3505
$ws.Range("J97").Value = 2950
$ws.Range("K97").Value = 3505
$ws.Range("L97").Value = 2950
$ws.Range("M97").Value = -3009
$ws.Range("N97").Value = -3942

$ws.Range("H107").Value = 1009
$ws.Range("I107").Value = 807.1429000000001
$ws.Range("K107").Value = 807.1429000000001
$ws.Range("M107").Value = 1112.8571

$ws.Range("H113").Value = 1874.8096
$ws.Range("I113").Value = 1874.0588
$ws.Range("J113").Value = 1878
$ws.Range("K113").Value = 1874.0588
$ws.Range("L113").Value = 1878
$ws.Range("M113").Value = 295.9412
$ws.Range("N113").Value = -6218

$ws.Range("H132").Value = 2245.5908
$ws.Range("I132").Value = 1690.0714
$ws.Range("J132").Value = 3217.75
$ws.Range("K132").Value = 5070.2142
$ws.Range("L132").Value = 9653.25
$ws.Range("M132").Value = -2540.2142
$ws.Range("N132").Value = -14713.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1021.125
$ws.Range("I16").Value = 1052.5333
$ws.Range("J16").Value = 550
$ws.Range("K16").Value = 1052.5333
$ws.Range("L16").Value = 550
$ws.Range("M16").Value = -882.5333000000001
$ws.Range("N16").Value = -890

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7078.778
$ws.Range("I81").Value = 17461
$ws.Range("J81").Value = 1887.6666
$ws.Range("K81").Value = 34922
$ws.Range("L81").Value = 3775.3332
$ws.Range("M81").Value = -33861
$ws.Range("N81").Value = -5897.3332

$ws.Range("H84").Value = 7078.778
$ws.Range("I84").Value = 17461
$ws.Range("J84").Value = 1887.6666
$ws.Range("K84").Value = 174610
$ws.Range("L84").Value = 18876.666
$ws.Range("M84").Value = -169306
$ws.Range("N84").Value = -29484.666

$ws.Range("H104").Value = 29723.334
$ws.Range("J104").Value = 29723.334
$ws.Range("L104").Value = 29723.334
$ws.Range("N104").Value = -36711.334
